$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "261.01"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.11%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.02"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.53%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.697"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.24%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06223"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.58%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.749"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8541"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.87%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9130"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.05%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1405"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.10%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04850"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-3.81%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07093"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.11%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03107"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.27%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.35%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001529"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.85%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006171"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.17%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006025"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.26%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.35%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.175"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.04%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.167"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.43%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1310"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.91%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.091"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.41%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04230"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.21%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001213"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.52%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.03%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.38%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03933"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.30%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1111"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.15%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004121"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.20%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.19%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01388"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-7.62%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.62%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.03%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1869"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "41.48%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.03%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.03%"
